# Bulk Upload Template - Employees sheet:
#  - F5 (Current Role) changes from "Java Developer" to "Solution Developer"
#  - F12 (Current Role) is cleared (value removed, style kept)
#  - Selection moves from G12 to F12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

$ws.Range("F5").Value = "Solution Developer"
$ws.Range("F12").Value = ""

[void]$ws.Range("F12").Select()
